$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("15").Insert()

$ws.Range("A15").Value = "code"
$ws.Range("B15").Value = "mpcb"
$ws.Range("A15").Font.Bold = $false
$ws.Range("A15").Font.Size = 12

$ws.Range("B16").Select()
